# Apply the "update conf and csv" edit to the Netflix workbook:
#  - Append a new header row (39) and a new data row (40) below the
#    existing "Year / Asia_sub / Sub_growth / Asia_revenue / ..." table,
#    introducing a new shared string "Net_sub_growth ".
#  - Update the sheet view (top-left cell / active selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New header row (row 39) ------------------------------------------------
$ws.Range("A39").Value = "Year "
$ws.Range("B39").Value = "Net_sub_growth "
$ws.Range("C39").Value = "Net_revenue_growth"
$ws.Range("D39").Value = "Movie_growth"
$ws.Range("E39").Value = "Audience_growth"

# Match the formatting used by the header row of the table directly above
# (row 34), which uses the "0.00_);[Red]\(0.00\)" number format on every cell
# (escaped parens so this reuses the workbook's existing numFmtId 176
# instead of minting a duplicate).
$ws.Range("A39:E39").NumberFormat = "0.00_);[Red]\(0.00\)"

# --- New data row (row 40) --------------------------------------------------
$ws.Range("A40").Value = 2020
$ws.Range("B40").Value = 64
$ws.Range("C40").Value = 62
$ws.Range("D40").Value = -73
$ws.Range("E40").Value = -74

$ws.Range("A40").NumberFormat = "0_);[Red]\(0\)"
$ws.Range("B40").NumberFormat = "0.00_);[Red]\(0.00\)"
$ws.Range("C40:E40").NumberFormat = "0.00_ "

# --- Sheet view updates ------------------------------------------------------
# Scroll the window so row 27 is at the top (best-effort — window scroll
# position isn't always persisted) and select the new active cell, matching
# the saved selection/topLeftCell in the target sheet view.
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I41").Select()
